$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Status column (K): remaining "Offen" rows are now closed ---
$ws.Range("K3").Value = "Geschlossen"
$ws.Range("K4").Value = "Geschlossen"
$ws.Range("K5").Value = "Geschlossen"
$ws.Range("K6").Value = "Geschlossen"

# --- "Last modified date" column (L): refresh the dates ---
# Row 2 keeps its original text date (09.04.2025), rows 3-6 get updated
# real date values reflecting the latest edits to the risk table.
$ws.Range("L3").Value = 45812   # 04.06.2025
$ws.Range("L4").Value = 45822   # 14.06.2025
$ws.Range("L5").Value = 45812   # 04.06.2025
$ws.Range("L6").Value = 45823   # 15.06.2025

$ws.Range("L3:L6").NumberFormat = "m/d/yy"

# --- Right-align the whole "Last modified date" column (header + body) ---
$ws.Range("L1:L6").HorizontalAlignment = -4152

# --- Widen column L so the right-aligned dates/header are fully visible ---
$ws.Columns("L").ColumnWidth = 16.666666666666668
